# "Updated to MATLAB 2023b" — the MOC nozzle design data on Sheet1 was
# regenerated, so the cached A2:B16 values change, and the sheet picks up
# an explicit column width for A:B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A2:B16 values (regenerated by the updated MATLAB script).
$newValues = @{
    2  = @(1.3760406834250329,  1.4407606868944756)
    3  = @(2.033274321331886,   1.635890667751646)
    4  = @(2.3907382043656482,  1.7337558547909908)
    5  = @(2.7372971564433803,  1.8207165240049576)
    6  = @(3.0882649800088,     1.9008493403852094)
    7  = @(3.4517297632007469,  1.9756986257760056)
    8  = @(3.8331534197834287,  2.0457822230274285)
    9  = @(4.236896125644944,   2.1110776111208867)
    10 = @(4.6668580954829482,  2.1712121334138463)
    11 = @(5.126800411866272,   2.2255429706600935)
    12 = @(5.6205306047007726,  2.273187971651319)
    13 = @(6.1520257814507477,  2.3130309116000372)
    14 = @(6.7255257280586287,  2.3437112142351899)
    15 = @(7.3456120296450624,  2.3636025870620045)
    16 = @(8.0172819828917472,  2.3707824297215985)
}

foreach ($row in $newValues.Keys) {
    $pair = $newValues[$row]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}

# Columns A:B get an explicit width of 7 characters.
$ws.Columns("A:B").ColumnWidth = 6.1666666666667
